# Apply updated frequency table values (rows 2-5, columns B:X) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(0.00681302043906132,0.00151400454201363,0.00151400454201363,0.000757002271006813,0.000757002271006813,0.00151400454201363,0.0416351249053747,0.0355791067373202,0.0105980317940954,0.0227100681302044,0.012112036336109,0.00454201362604088,0.00529901589704769,0.994700984102952,0.000757002271006813,0,0,0,0,0.00151400454201363,0.934140802422407,0.0151400454201363,0.0060560181680545),
    @(0.987887963663891,0.00757002271006813,0.0060560181680545,0.00681302043906132,0.989401968205905,0.98107494322483,0.0060560181680545,0.945495836487509,0.0280090840272521,0.000757002271006813,0,0.000757002271006813,0,0,0.998485995457986,0,0.0060560181680545,0.00151400454201363,0.999242997728993,0.993186979560939,0,0.00908402725208176,0),
    @(0.000757002271006813,0.00529901589704769,0,0.00227100681302044,0.000757002271006813,0.00529901589704769,0.94776684330053,0.00302800908402725,0.000757002271006813,0.9666919000757,0.987130961392884,0.993943981831946,0.993943981831946,0.000757002271006813,0.000757002271006813,0,0,0,0,0.000757002271006813,0.0643451930355791,0.9666919000757,0.98107494322483),
    @(0.00378501135503407,0.985616956850871,0.992429977289932,0.990158970476911,0.00908402725208176,0.012112036336109,0.00454201362604088,0.0158970476911431,0.960635881907646,0.00908402725208176,0,0,0.000757002271006813,0.00454201362604088,0,1,0.993943981831946,0.998485995457986,0.000757002271006813,0.00454201362604088,0.00151400454201363,0.00908402725208176,0.0128690386071158)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowIndex = $r + 2
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $colIndex = $c + 2
        $ws.Cells.Item($rowIndex, $colIndex).Value = $rowValues[$c]
    }
}

$wb.Save()
